$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (existing D:K shift to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting into the new D:E columns from column F (matches row style: date fmt for header rows, number fmt otherwise)
$ws.Range("F7:F7").Copy()
$ws.Range("D7:E7").PasteSpecial(-4122)
$ws.Range("F8:F35").Copy()
$ws.Range("D8:E35").PasteSpecial(-4122)
$ws.Range("F38:F38").Copy()
$ws.Range("D38:E38").PasteSpecial(-4122)
$ws.Range("F39:F77").Copy()
$ws.Range("D39:E77").PasteSpecial(-4122)
$ws.Range("F80:F80").Copy()
$ws.Range("D80:E80").PasteSpecial(-4122)
$ws.Range("F81:F102").Copy()
$ws.Range("D81:E102").PasteSpecial(-4122)

# Populate values for the new D and E columns
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 10600
$ws.Range("E8").Value = 9100
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = 24800
$ws.Range("E12").Value = 23700
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 32300
$ws.Range("E17").Value = 31300
$ws.Range("D18").Value = -21700
$ws.Range("E18").Value = -22200
$ws.Range("D20").Value = 400
$ws.Range("E20").Value = 200
$ws.Range("D21").Value = -19300
$ws.Range("E21").Value = -20000
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("D23").Value = -21300
$ws.Range("E23").Value = -21900
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -21300
$ws.Range("E26").Value = -21900
$ws.Range("D27").Value = -21300
$ws.Range("E27").Value = -21900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -400
$ws.Range("E32").Value = -200
$ws.Range("D33").Value = -21300
$ws.Range("E33").Value = -21900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -21300
$ws.Range("E35").Value = -21900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 85800
$ws.Range("E41").Value = 55700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 17200
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 6800
$ws.Range("E45").Value = 6600
$ws.Range("D46").Value = 92700
$ws.Range("E46").Value = 79400
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("D48").Value = 26300
$ws.Range("E48").Value = 28100
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 1500
$ws.Range("E52").Value = 1500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 120500
$ws.Range("E54").Value = 109000
$ws.Range("D57").Value = 6400
$ws.Range("E57").Value = 6500
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 35600
$ws.Range("E59").Value = 32500
$ws.Range("D60").Value = 42000
$ws.Range("E60").Value = 39000
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 126500
$ws.Range("E62").Value = 100600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 168500
$ws.Range("E66").Value = 139600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -389400
$ws.Range("E72").Value = -368100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = -48000
$ws.Range("E76").Value = -30600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -21300
$ws.Range("E81").Value = -21900
$ws.Range("D83").Value = 2000
$ws.Range("E83").Value = 2000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 12700
$ws.Range("E89").Value = -23100
$ws.Range("D91").Value = 200
$ws.Range("E91").Value = -500
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 17400
$ws.Range("E94").Value = 20500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 300
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 30100
$ws.Range("E102").Value = -2200

# Data corrections for row 91 (not a pure shift of the old G/H values)
$ws.Range("I91").Value = -500
$ws.Range("J91").Value = -1300
